# FOSWEC Sign Conventions - slide 3: add flap/motor sign annotations
#
# Adds two "+flap.*Flap_deg" labels, two red arc/arrow indicators, two
# "+motor.*Pos_rad" labels, and a "Cycling low power resets position zero"
# note to the 3rd slide, matching the variable-renaming / sign-convention
# commit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# EMU-per-point constant used by the PowerPoint object model (Left/Top/
# Width/Height and AddShape/AddTextbox all traffic in points, while the
# underlying OOXML stores English Metric Units).
$emuPerPt = 12700.0

function ToPt([double]$emu) {
    return $emu / $emuPerPt
}

# ---------------------------------------------------------------------
# 1) "+flap.bowFlap_deg" textbox
# ---------------------------------------------------------------------
$tb1 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$tb1.Name = "TextBox 1"
$tb1.TextFrame.WordWrap = $false
$tb1.TextFrame.AutoSize = 1
$tb1.Fill.Visible = $false
$tr1 = $tb1.TextFrame.TextRange
$tr1.Text = "+"
$tr1.InsertAfter("flap.bowFlap_deg") | Out-Null
$tb1.Left = ToPt(3894780)
$tb1.Top = ToPt(3244333)
$tb1.Width = ToPt(1974323)
$tb1.Height = ToPt(369332)

# ---------------------------------------------------------------------
# 2) "+flap.aftFlap_deg" textbox
# ---------------------------------------------------------------------
$tb2 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$tb2.Name = "TextBox 13"
$tb2.TextFrame.WordWrap = $false
$tb2.TextFrame.AutoSize = 1
$tb2.Fill.Visible = $false
$tr2 = $tb2.TextFrame.TextRange
$tr2.Text = "+"
$tr2.InsertAfter("flap.aftFlap_deg") | Out-Null
$tb2.Left = ToPt(9800948)
$tb2.Top = ToPt(3316778)
$tb2.Width = ToPt(1823191)
$tb2.Height = ToPt(369332)

# ---------------------------------------------------------------------
# 3) Arc 14 - red sign-convention arc (duplicate of the existing purple
#    "Arc 11" so the adj1/adj2 guides, p:style and line-end markers all
#    come along for free)
# ---------------------------------------------------------------------
$arcSrc1 = $s.Shapes.Item("Arc 11")
$arc1Range = $arcSrc1.Duplicate()
$arc1 = $arc1Range.Item(1)
$arc1.Name = "Arc 14"
$arc1.Left = ToPt(6241233)
$arc1.Top = ToPt(4653547)
$arc1.Width = ToPt(1573626)
$arc1.Height = ToPt(1573626)
$arc1.Line.ForeColor.RGB = 255

# ---------------------------------------------------------------------
# 4) Arc 15 - red sign-convention arc
# ---------------------------------------------------------------------
$arcSrc2 = $s.Shapes.Item("Arc 11")
$arc2Range = $arcSrc2.Duplicate()
$arc2 = $arc2Range.Item(1)
$arc2.Name = "Arc 15"
$arc2.Left = ToPt(4667607)
$arc2.Top = ToPt(4618768)
$arc2.Width = ToPt(1573626)
$arc2.Height = ToPt(1573626)
$arc2.Line.ForeColor.RGB = 255

# ---------------------------------------------------------------------
# 5) "+motor.bowPos_rad" textbox
# ---------------------------------------------------------------------
$tb3 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$tb3.Name = "TextBox 2"
$tb3.TextFrame.WordWrap = $false
$tb3.TextFrame.AutoSize = 1
$tb3.Fill.Visible = $false
$tr3 = $tb3.TextFrame.TextRange
$tr3.Text = "+"
$tr3.InsertAfter("motor.bowPos_rad") | Out-Null
$tb3.Left = ToPt(3994171)
$tb3.Top = ToPt(4265317)
$tb3.Width = ToPt(2074158)
$tb3.Height = ToPt(369332)

# ---------------------------------------------------------------------
# 6) "+motor.aftPos_rad" textbox
# ---------------------------------------------------------------------
$tb4 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$tb4.Name = "TextBox 18"
$tb4.TextFrame.WordWrap = $false
$tb4.TextFrame.AutoSize = 1
$tb4.Fill.Visible = $false
$tr4 = $tb4.TextFrame.TextRange
$tr4.Text = "+"
$tr4.InsertAfter("motor.aftPos_rad") | Out-Null
$tb4.Left = ToPt(6096000)
$tb4.Top = ToPt(4259994)
$tb4.Width = ToPt(1923027)
$tb4.Height = ToPt(369332)

# ---------------------------------------------------------------------
# 7) "Cycling low power resets position zero" note
# ---------------------------------------------------------------------
$tb5 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$tb5.Name = "TextBox 3"
$tb5.TextFrame.WordWrap = $false
$tb5.TextFrame.AutoSize = 1
$tb5.Fill.Visible = $false
$tb5.TextFrame.TextRange.Text = "Cycling low power resets position zero"
$tb5.Left = ToPt(3894780)
$tb5.Top = ToPt(6432672)
$tb5.Width = ToPt(3771482)
$tb5.Height = ToPt(369332)
